$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1438.3125
$ws.Range("J32").Value = 1655.1666
$ws.Range("L32").Value = 1655.1666
$ws.Range("N32").Value = -2307.1666
$ws.Range("H33").Value = 141.625
$ws.Range("I33").Value = 151.21428
$ws.Range("J33").Value = 74.5
$ws.Range("K33").Value = 151.21428
$ws.Range("L33").Value = 74.5
$ws.Range("M33").Value = 77.78572
$ws.Range("N33").Value = -532.5
$ws.Range("H40").Value = 5224.136
$ws.Range("I40").Value = 3590.6365
$ws.Range("K40").Value = 3590.6365
$ws.Range("M40").Value = -3415.6365
$ws.Range("H41").Value = 15153155
$ws.Range("I41").Value = 643.55554
$ws.Range("K41").Value = 643.55554
$ws.Range("M41").Value = -203.55554
$ws.Range("H86").Value = 2292207.8
$ws.Range("I86").Value = 3628.3333
$ws.Range("J86").Value = 3763437.2
$ws.Range("K86").Value = 3628.3333
$ws.Range("L86").Value = 3763437.2
$ws.Range("M86").Value = -2505.3333
$ws.Range("N86").Value = -3765683.2
$ws.Range("H89").Value = 2292207.8
$ws.Range("I89").Value = 3628.3333
$ws.Range("J89").Value = 3763437.2
$ws.Range("K89").Value = 18141.6665
$ws.Range("L89").Value = 18817186
$ws.Range("M89").Value = -12525.6665
$ws.Range("N89").Value = -18828418
$ws.Range("H92").Value = 330.2
$ws.Range("I92").Value = 266.875
$ws.Range("J92").Value = 402.57144
$ws.Range("K92").Value = 266.875
$ws.Range("L92").Value = 402.57144
$ws.Range("M92").Value = 981.125
$ws.Range("N92").Value = -2898.57144
$ws.Range("H127").Value = 3490
$ws.Range("I127").Value = 1392.6
$ws.Range("K127").Value = 4177.799999999999
$ws.Range("M127").Value = 782.2000000000007
$ws.Range("H129").Value = 2165.9375
$ws.Range("I129").Value = 1236.7142
$ws.Range("J129").Value = 2888.6667
$ws.Range("K129").Value = 3710.1426
$ws.Range("L129").Value = 8666.000100000001
$ws.Range("M129").Value = 1289.8574
$ws.Range("N129").Value = -18666.0001
$ws.Range("H132").Value = 4554.8
$ws.Range("I132").Value = 4593.7026
$ws.Range("K132").Value = 13781.1078
$ws.Range("M132").Value = -11251.1078
$ws.Range("H138").Value = 4835.3896
$ws.Range("I138").Value = 2172.0476
$ws.Range("J138").Value = 6307.237
$ws.Range("K138").Value = 6516.1428
$ws.Range("L138").Value = 18921.711
$ws.Range("M138").Value = -1376.1428
$ws.Range("N138").Value = -29201.711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1476.9661
$ws.Range("I32").Value = 1519.0364
$ws.Range("K32").Value = 1519.0364
$ws.Range("M32").Value = -1232.0364
$ws.Range("H45").Value = 36218.31
$ws.Range("I45").Value = 49130.094
$ws.Range("J45").Value = 2324.875
$ws.Range("K45").Value = 49130.094
$ws.Range("L45").Value = 2324.875
$ws.Range("M45").Value = -48753.094
$ws.Range("N45").Value = -3078.875
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1960.3334
$ws.Range("I20").Value = 1306.5714
$ws.Range("J20").Value = 4248.5
$ws.Range("K20").Value = 1306.5714
$ws.Range("L20").Value = 4248.5
$ws.Range("M20").Value = -1059.5714
$ws.Range("N20").Value = -4742.5
$ws.Range("H80").Value = 983.3077
$ws.Range("I80").Value = 1409.5
$ws.Range("J80").Value = 618
$ws.Range("K80").Value = 1409.5
$ws.Range("L80").Value = 618
$ws.Range("M80").Value = -411.5
$ws.Range("N80").Value = -2614
$ws.Range("H83").Value = 983.3077
$ws.Range("I83").Value = 1409.5
$ws.Range("J83").Value = 618
$ws.Range("K83").Value = 7047.5
$ws.Range("L83").Value = 3090
$ws.Range("M83").Value = -2055.5
$ws.Range("N83").Value = -13074
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H134").Value = 33903
$ws.Range("I134").Value = 1638.1154
$ws.Range("K134").Value = 4914.3462
$ws.Range("M134").Value = -2379.3462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2624.1667
$ws.Range("I16").Value = 2213.1428
$ws.Range("J16").Value = 3199.6
$ws.Range("K16").Value = 2213.1428
$ws.Range("L16").Value = 3199.6
$ws.Range("M16").Value = -1926.1428
$ws.Range("N16").Value = -3773.6
$ws.Range("H52").Value = 88849.5
$ws.Range("H58").Value = 7862.533
$ws.Range("I58").Value = 2996.8462
$ws.Range("J58").Value = 11583.353
$ws.Range("K58").Value = 2996.8462
$ws.Range("L58").Value = 11583.353
$ws.Range("M58").Value = -2793.8462
$ws.Range("N58").Value = -11989.353
$ws.Range("H99").Value = 6824.5
$ws.Range("I99").Value = 5332.3335
$ws.Range("K99").Value = 5332.3335
$ws.Range("M99").Value = -3834.3335
$ws.Range("H107").Value = 1198.28
$ws.Range("I107").Value = 753.6111
$ws.Range("K107").Value = 753.6111
$ws.Range("M107").Value = 1166.3889
$ws.Range("H113").Value = 2624.1667
$ws.Range("I113").Value = 2213.1428
$ws.Range("J113").Value = 3199.6
$ws.Range("K113").Value = 2213.1428
$ws.Range("L113").Value = 3199.6
$ws.Range("M113").Value = -43.14280000000008
$ws.Range("N113").Value = -7539.6
$ws.Range("H122").Value = 3020.0667
$ws.Range("I122").Value = 2600.5
$ws.Range("J122").Value = 3859.2
$ws.Range("K122").Value = 7801.5
$ws.Range("L122").Value = 11577.6
$ws.Range("M122").Value = -5351.5
$ws.Range("N122").Value = -16477.6
$ws.Range("H126").Value = 6824.5
$ws.Range("I126").Value = 5332.3335
$ws.Range("K126").Value = 15997.0005
$ws.Range("M126").Value = -13527.0005
$ws.Range("H132").Value = 3165.4075
$ws.Range("I132").Value = 1792.8334
$ws.Range("K132").Value = 5378.5002
$ws.Range("M132").Value = -2848.5002
$ws.Range("H136").Value = 7862.533
$ws.Range("I136").Value = 2996.8462
$ws.Range("J136").Value = 11583.353
$ws.Range("K136").Value = 8990.5386
$ws.Range("L136").Value = 34750.05899999999
$ws.Range("M136").Value = -6440.5386
$ws.Range("N136").Value = -39850.05899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 253055.25
$ws.Range("I116").Value = 253055.25
$ws.Range("K116").Value = 759165.75
$ws.Range("M116").Value = -755723.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1255316.6
$ws.Range("I80").Value = 1003003.9
$ws.Range("K80").Value = 1003003.9
$ws.Range("M80").Value = -1002005.9
$ws.Range("H83").Value = 1255316.6
$ws.Range("I83").Value = 1003003.9
$ws.Range("K83").Value = 5015019.5
$ws.Range("M83").Value = -5010027.5
$ws.Range("H97").Value = 549.26666
$ws.Range("J97").Value = 833.75
$ws.Range("L97").Value = 833.75
$ws.Range("N97").Value = -1825.75
$ws.Range("H122").Value = 584184.8
$ws.Range("I122").Value = 739648.8
$ws.Range("J122").Value = 1194.75
$ws.Range("K122").Value = 2218946.4
$ws.Range("L122").Value = 3584.25
$ws.Range("M122").Value = -2216496.4
$ws.Range("N122").Value = -8484.25
$ws.Range("H126").Value = 9529
$ws.Range("J126").Value = 11693.154
$ws.Range("L126").Value = 35079.462
$ws.Range("N126").Value = -40019.462
$ws.Range("H132").Value = 100380.43
$ws.Range("I132").Value = 29610.5
$ws.Range("J132").Value = 525000
$ws.Range("K132").Value = 88831.5
$ws.Range("L132").Value = 1575000
$ws.Range("M132").Value = -86301.5
$ws.Range("N132").Value = -1580060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2841.9714
$ws.Range("I46").Value = 2294.5
$ws.Range("J46").Value = 3768.4614
$ws.Range("K46").Value = 2294.5
$ws.Range("L46").Value = 3768.4614
$ws.Range("M46").Value = -2106.5
$ws.Range("N46").Value = -4144.4614
$ws.Range("H55").Value = 33333860
$ws.Range("J55").Value = 83334310
$ws.Range("L55").Value = 83334310
$ws.Range("N55").Value = -83334656
$ws.Range("H100").Value = 80426.664
$ws.Range("I100").Value = 124322.22
$ws.Range("K100").Value = 124322.22
$ws.Range("M100").Value = -123781.22
$ws.Range("H132").Value = 6442.4585
$ws.Range("I132").Value = 5472.3076
$ws.Range("J132").Value = 7589
$ws.Range("K132").Value = 16416.9228
$ws.Range("L132").Value = 22767
$ws.Range("M132").Value = -13886.9228
$ws.Range("N132").Value = -27827

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1293.3077
$ws.Range("I100").Value = 1355.4546
$ws.Range("J100").Value = 951.5
$ws.Range("K100").Value = 2710.9092
$ws.Range("L100").Value = 1903
$ws.Range("M100").Value = -2169.9092
$ws.Range("N100").Value = -2985
$ws.Range("H136").Value = 321868.78
$ws.Range("I136").Value = 325703.9
$ws.Range("K136").Value = 977111.7000000001
$ws.Range("M136").Value = -974561.7000000001
